$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.748.69"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.076.49"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.95"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.19"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("D12").Value = "2.373.20"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.44"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.95"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "2.080.55"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "37.705.07"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.62"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.71"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.99"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("E27").Value = "  +10.62%  "
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0994"
$ws.Range("E40").Value = "  +6.76%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.62"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D44").Value = "1.451.63"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.26"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.63"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.03"
$ws.Range("E51").Value = "  +7.37%  "
